$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 10 ("Are we over-fitting the data?") - Title 1
# Split the title into styled runs and append a new sentence.
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$titleShape = $s10.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange

# Reset first so the engine doesn't keep an internal run-boundary "ghost"
# at the length of the old text when we immediately re-split the new text
# into runs below (avoids spurious extra run splits on save).
$titleRange.Text = "x"
$titleRange.Text = "Are we over-fitting the data? (Bayesian approach has another criteria)"

$titleRange.Characters(8, 12).Font.Color.RGB = 255
$titleRange.Characters(31, 40).Font.Size = 14

# ---------------------------------------------------------------------------
# Slide 11 ("Codes you might want to use sometime") - Content Placeholder 3
# ---------------------------------------------------------------------------
$s11 = $p.Slides.Item(11)
$bodyShape = $s11.Shapes.Item(2)
$bodyRange = $bodyShape.TextFrame.TextRange

$para1 = $bodyRange.Paragraphs(1)
$para1.Text = "x"
$para1.Text = "If your likelihood is a mixture of multiple Gaussians: (approximate any distribution by multiple Gaussians)"

$p1start = $para1.Start
$bodyRange.Characters($p1start + 35, 18).Font.Color.RGB = 255
$bodyRange.Characters($p1start + 56, 50).Font.Size = 12

$para5 = $bodyRange.Paragraphs(5)
$para5.Text = "x"
$para5.Text = "Bootstrap and jackknife methods: see code in section 4.5"

$p5start = $para5.Start
$bodyRange.Characters($p5start, 31).Font.Color.RGB = 255
